$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-393). The sheet was refreshed and the stored "last changed"
# date moved forward by one day, from serial 45188 (2023-09-19) to
# 45189 (2023-09-20), for every single row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Range("C2:C$lastRow").Value = 45189
